{"js": "// refactor marine v inland logic with Boolean and change one style font size\n//\n// 1) The stray \"_GoBack\" bookmark that Word leaves at the last edit\n//    position (around {{table:priorities}}) is removed.\n// 2) A fresh \"_GoBack\" bookmark is dropped at the new last-edit spot,\n//    splitting the run \"... measured the current condition ...\" right\n//    after the word \"current\".\n// 3) The HyperlinkSource paragraph style picks up an explicit 11pt\n//    (sz=22) font size.\n\n// --- 1. Remove the old \"_GoBack\" bookmark -------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Re-insert \"_GoBack\" at the new position --------------------------\n// Locate the unique run of text and split right after \"current\".\nconst searchResults = context.document.body.search(\n  \"This report card measured the current\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const found = searchResults.items[0];\n  const splitPoint = found.getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 3. Bump the HyperlinkSource style font size to 11pt (sz=22) --------\nconst hyperlinkSourceStyle = context.document.getStyles().getByNameOrNullObject(\"HyperlinkSource\");\nhyperlinkSourceStyle.load(\"nameLocal\");\nawait context.sync();\n\nif (!hyperlinkSourceStyle.isNullObject) {\n  hyperlinkSourceStyle.font.size = 11;\n  await context.sync();\n}\n", "ps1": "# refactor marine v inland logic with Boolean and change one style font size\n#\n# 1) The stray \"_GoBack\" bookmark that Word leaves at the last edit\n#    position (around {{table:priorities}}) is removed.\n# 2) A fresh \"_GoBack\" bookmark is dropped at the new last-edit spot,\n#    splitting the run \"... measured the current condition ...\" right\n#    after the word \"current\".\n# 3) The HyperlinkSource paragraph style picks up an explicit 11pt\n#    (sz=22) font size.\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the old \"_GoBack\" bookmark --------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2. Re-insert \"_GoBack\" at the new position --------------------------\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"This report card measured the current\")\n\nif ($found) {\n  $insertionPoint = $d.Range($searchRange.End, $searchRange.End)\n  $d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n\n# --- 3. Bump the HyperlinkSource style font size to 11pt (sz=22) --------\n$hyperlinkSourceStyle = $d.Styles(\"HyperlinkSource\")\n$hyperlinkSourceStyle.Font.Size = 11\n"}
